# Skill.xlsx: "finish the crt dam skill" (close #38)
# - add a new 暴击 (crit) skill row (Id 55900031) right before the 5599xxxx block,
#   shifting all following rows down by one.
# - grow the worksheet Table to cover the new row.
# - nudge the view/selection to the newly added row.
# - mark the workbook window minimized (best effort).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# --- insert a new row 149, pushing 149..164 down to 150..165 ---------------
$ws.Rows.Item(149).Insert()

# Clone formatting from row 148 (the row that used to sit directly above the
# new one) so the new row matches its neighbours' look (borders/fonts/etc).
$ws.Range("A148:Z148").Copy()
$ws.Range("A149:Z149").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- fill in the new row's data --------------------------------------------
# Set the brand-new text values first, in left-to-right / top-to-bottom
# column order, so new shared-string entries line up the same way Excel
# would naturally append them while the row is being authored.
$ws.Cells.Item(149, 8).Value  = "s.CrtDamAddRate+=0.5;"    # H: OnAdd
$ws.Cells.Item(149, 2).Value  = "暴击"                      # B: Name
$ws.Cells.Item(149, 19).Value = "暴击伤害提升50%"            # S: GetDescript
$ws.Cells.Item(149, 25).Value = "baoji"                    # Y: Icon

$ws.Cells.Item(149, 1).Value  = 55900031                   # A: Id
$ws.Cells.Item(149, 3).Value  = "特殊"                      # C: Type
$ws.Cells.Item(149, 17).Value = "Active"                   # Q: Active
# Leading apostrophe forces literal text "true" instead of a Boolean, to
# match the rest of the IsRandom column (stored as the text "true").
$ws.Cells.Item(149, 18).Value = "'true"                    # R: IsRandom
$ws.Cells.Item(149, 24).Value = 5                           # X: Mark

# --- grow the listobject/table to include the new row ----------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:Z165"))

# --- update worksheet dimension / selection to match the new layout --------
$ws.Range("A149").Select()

# --- reflect the view being minimized when the file was saved (best effort)
$excel.ActiveWindow.WindowState = -4140   # xlMinimized
